# dsa tries and bits
# Add two new rows (72, 73) to the LeetCode tracking sheet:
#   1268. Search Suggestions System                 (Tries)
#   1318. Minimum Flips to Make a OR b Equal to c    (Bit Manipulation)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 71) down into the
# two new rows first, so Difficulty fill / hyperlink font come along.
$ws.Range("A71:E71").Copy()
$ws.Range("A72:E73").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 72: 1268. Search Suggestions System -----------------------------
$ws.Range("A72").Value = "1268. Search Suggestions System"
$ws.Range("B72").Value = "Medium"
$ws.Range("C72").Value = "Tries"
$ws.Range("E72").Value = "https://leetcode.com/problems/search-suggestions-system/solutions/436151/java-python-3-simple-trie-and-binary-search-codes-w-comment-and-brief-analysis/?envType=study-plan-v2&envId=leetcode-75 "
$ws.Range("D72").Value = "Classic Trie, but can be solved with Binary Search. Use DS principles. Class Trie with linked list for suggestions."

# --- Row 73: 1318. Minimum Flips to Make a OR b Equal to c ----------------
$ws.Range("A73").Value = "1318. Minimum Flips to Make a OR b Equal to c"
$ws.Range("B73").Value = "Medium"
$ws.Range("C73").Value = "Bit Manipulation"
$ws.Range("D73").Value = "Quick trick is use Integer.bitCount(), but from first principles, use a bit mask."
$ws.Range("E73").Value = "https://leetcode.com/problems/minimum-flips-to-make-a-or-b-equal-to-c/solutions/477690/java-python-3-bit-manipulation-w-explanation-and-analysis/?envType=study-plan-v2&envId=leetcode-75 "

# --- Turn E72 / E73 into real hyperlinks (this resets their style, so we
#     re-apply the copied Hyperlink-column formatting afterwards) ---------
$ws.Hyperlinks.Add($ws.Range("E72"), "https://leetcode.com/problems/search-suggestions-system/solutions/436151/java-python-3-simple-trie-and-binary-search-codes-w-comment-and-brief-analysis/?envType=study-plan-v2&envId=leetcode-75") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E73"), "https://leetcode.com/problems/minimum-flips-to-make-a-or-b-equal-to-c/solutions/477690/java-python-3-bit-manipulation-w-explanation-and-analysis/?envType=study-plan-v2&envId=leetcode-75") | Out-Null

$ws.Range("E72").Value = "https://leetcode.com/problems/search-suggestions-system/solutions/436151/java-python-3-simple-trie-and-binary-search-codes-w-comment-and-brief-analysis/?envType=study-plan-v2&envId=leetcode-75 "
$ws.Range("E73").Value = "https://leetcode.com/problems/minimum-flips-to-make-a-or-b-equal-to-c/solutions/477690/java-python-3-bit-manipulation-w-explanation-and-analysis/?envType=study-plan-v2&envId=leetcode-75 "

$ws.Range("E71").Copy()
$ws.Range("E72:E73").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the selection the same way the original author's session ended up.
$ws.Range("E79").Select() | Out-Null
